# Continuous Lab support: drop the "key" (win2012ssh) column from the
# UserList table. The sheet goes from 4 columns (email, key, role,
# awsAccount) to 3 (email, role, awsAccount); column C (role) and column D
# (awsAccount) shift left into B and C respectively, carrying their values,
# shared-string type and cell styles (incl. the text-format style on the
# awsAccount column) along with them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("B").Delete() | Out-Null

# Restore the active-cell selection recorded for this sheet after the edit.
$ws.Range("N14").Select() | Out-Null
